$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "베이지안 최적화에 기반한 HyperOpt를 활용한 하이퍼 파라미터 튜닝"
$ws.Range("E4").Value = "https://teddylee777.github.io/thoughts/hyper-opt"

$ws.Range("D16").Value = "Sanity checks for saliency maps 내용정리 [XAI-6 (2)]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/148"

$ws.Range("D28").Value = "EfficientNet : Rethinking Model Scaling for Convolutional Neural Networks 논문 리뷰"
$ws.Range("E28").Value = "https://ropiens.tistory.com/110"

$ws.Range("D29").Value = "[만화] 인턴일기 8~12"
$ws.Range("E29").Value = "https://blog.promedius.ai/intern-life-2/"

$ws.Range("D46").Value = "[씨젠] 2021년 05월, 생물정보학(Bioinformatics 채용), Bioinformatics Data Engineer, Bioinformatics SW Engineer"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/392"

$ws.Range("D51").Value = "[python] pandas dataframe 모든 열 또는 행 보여주기"
$ws.Range("E51").Value = "https://bskyvision.com/1176"
